# edit.ps1
# Implements commit "feat: add 2022-Q4 data":
#   - Insert a new "2022-Q4" sheet (with its own fund table) right before "2022-Q3".
#   - Insert a new summary row for "2022-Q4" at the top of the "总计" (totals) sheet,
#     pushing the existing quarters down by one row.

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force a numeric-looking string to be stored as TEXT (keeps leading zeros /
    # trailing decimal zeros) the same way a user typing  '4.70  into Excel would.
    $Cell.Value = "'" + $Text
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q4" worksheet by duplicating the "2022-Q3" sheet
#    (this keeps header / column styling identical) and then overwriting
#    its data with the 2022-Q4 numbers, trimming the extra row.
# ---------------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($sheetQ3)
# NOTE: after Copy(), $sheetQ3 now refers to the newly-inserted duplicate sheet
# (the object reference tracks the sheet that occupies that position), so we
# rename it directly and look the original "2022-Q3" sheet back up by name.
$newSheet = $sheetQ3
$newSheet.Name = "2022-Q4"
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")

# The copied sheet has 9 data rows (rows 2-10); 2022-Q4 only needs 8 (rows 2-9).
$newSheet.Rows.Item(10).Delete()

# Header row (row 1) text stays the same (基金代码, 基金名称, ...) - leave as-is.

$q4Data = @(
    @("002376", "国寿安保核心产业灵活配置混合", "4.70", "88.55", "4.65", "0.2186", 2),
    @("005175", "国寿安保消费新蓝海灵活配置混合", "0.70", "89.58", "5.96", "0.0417", 2),
    @("004258", "国寿安保稳嘉混合A", "2.16", "23.32", "1.07", "0.0231", 8),
    @("006813", "博时汇悦回报混合", "0.73", "75.26", "2.63", "0.0192", 8),
    @("004301", "国寿安保稳信混合A", "1.48", "22.04", "1.01", "0.0149", 6),
    @("004259", "国寿安保稳嘉混合C", "0.01", "23.32", "1.07", "0.0001", 8),
    @("004302", "国寿安保稳信混合C", "0.01", "22.04", "1.01", "0.0001", 6),
    @("015406", "国寿安保稳信混合E", "0.00", "22.04", "1.01", 0, 6)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $rec = $q4Data[$i]

    $newSheet.Cells.Item($row, 1).Value = $i          # A: index number (0-based)
    Set-TextValue $newSheet.Cells.Item($row, 2) $rec[0]   # B: 基金代码 (fund code)
    $newSheet.Cells.Item($row, 3).Value = $rec[1]         # C: 基金名称 (fund name)
    Set-TextValue $newSheet.Cells.Item($row, 4) $rec[2]   # D: 基金规模
    Set-TextValue $newSheet.Cells.Item($row, 5) $rec[3]   # E: 股票总仓位
    Set-TextValue $newSheet.Cells.Item($row, 6) $rec[4]   # F: 仓位占比

    if ($row -eq 9) {
        $newSheet.Cells.Item($row, 7).Value = $rec[5]     # G: last row is a real 0 number
    } else {
        Set-TextValue $newSheet.Cells.Item($row, 7) $rec[5]
    }
    $newSheet.Cells.Item($row, 8).Value = $rec[6]         # H: 仓位排名 (number)
}

# ---------------------------------------------------------------------------
# 2) Insert a new "2022-Q4" row at the top of the 总计 (totals) sheet, pushing
#    the other quarters down by one row, and fill it with the new totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# NOTE: column A is just a fixed 0-based row-index (A2=0, A3=1, ... A8=6) and
# does NOT travel with the shifted data - only B/C/D (quarter, count, value)
# move down one row.
$lastRow = 7   # existing data occupies rows 2-7
for ($r = $lastRow; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1

    # Copy formatting from the source row first (so the newly-created row 8
    # below the old data keeps the same look), then stamp in the values.
    $total.Range("A" + $src).Copy()
    $total.Range("A" + $dst).PasteSpecial(-4122)

    $bVal = $total.Cells.Item($src, 2).Value2
    $cVal = $total.Cells.Item($src, 3).Value2
    $dVal = $total.Cells.Item($src, 4).Value2

    $total.Cells.Item($dst, 1).Value = $dst - 2
    $total.Cells.Item($dst, 2).Value = $bVal
    $total.Cells.Item($dst, 3).Value = $cVal
    $total.Cells.Item($dst, 4).Value = $dVal
}

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.32
